$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header labels for the team record columns: AD=Wins, AE=Losses, AF=Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the existing header formatting (bold font + border + centered alignment)
# from the neighboring header cell (AC1) rather than building a brand-new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the team's Wins/Losses/Ties for every player row (2 through 51)
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 68
    $ws.Cells.Item($row, 31).Value = 94
    $ws.Cells.Item($row, 32).Value = 0
}
